# Apply the data corrections to column F (dSF) as described in the commit:
# "repull data, push all data, mean calculation"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -4
$ws.Range("F3").Value = -2
$ws.Range("F6").Value = -4
$ws.Range("F9").Value = 2
$ws.Range("F11").Value = -2
